# Update "currentforecastonactivity" report sheet:
#  - refresh report date (B1)
#  - remove the "NOT" activity row entirely (rows below shift up)
#  - refresh activity counts (col B) and hour index (col C) for every row,
#    aligning the "hour" column to the new index (14 -> 16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "NOT" row (row 8) - all following rows shift up one position.
$ws.Rows("8").Delete()

# Refresh the report date in the header row. Force text so Excel doesn't
# reinterpret the dd/mm/yyyy-looking string as a serial date value.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "07/04/2023"

# Refresh activity totals (col B) and the aligned hour index (col C).
$ws.Range("B2").Value = 2139
$ws.Range("C2").Value = 16

$ws.Range("B3").Value = 181
$ws.Range("C3").Value = 16

$ws.Range("B4").Value = 461
$ws.Range("C4").Value = 16

$ws.Range("B5").Value = 1302
$ws.Range("C5").Value = 16

$ws.Range("B6").Value = 1109
$ws.Range("C6").Value = 16

$ws.Range("B7").Value = 289
$ws.Range("C7").Value = 16

$ws.Range("B8").Value = 905
$ws.Range("C8").Value = 16

$ws.Range("B9").Value = 141
$ws.Range("C9").Value = 16

$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 16

$ws.Range("B11").Value = 96
$ws.Range("C11").Value = 16
